$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "2023_01_02"

$ws.Range("D2").Value = 5.78
$ws.Range("E2").Value = 2.89
$ws.Range("D3").Value = 2.66
$ws.Range("E3").Value = 3.75
$ws.Range("D4").Value = 8.199999999999999
$ws.Range("E4").Value = 5.68
$ws.Range("D5").Value = 13
$ws.Range("E5").Value = 2.07
$ws.Range("D6").Value = 11.73
$ws.Range("E6").Value = 2.19
$ws.Range("D7").Value = 7.67
$ws.Range("E7").Value = 5.49
$ws.Range("D8").Value = 4.7
$ws.Range("E8").Value = 1.45
$ws.Range("D9").Value = 15.48
$ws.Range("E9").Value = 1.99
$ws.Range("D10").Value = 14.32
$ws.Range("E10").Value = 2.2
$ws.Range("D11").Value = 3.38
$ws.Range("D12").Value = 18
$ws.Range("E12").Value = 10.95
$ws.Range("E13").Value = 3.99
$ws.Range("D14").Value = 11.52
$ws.Range("E14").Value = 6.49
$ws.Range("D15").Value = 5.57
$ws.Range("E15").Value = 3.95
$ws.Range("D16").Value = 10
$ws.Range("D17").Value = 20.33
$ws.Range("E17").Value = 2.49
$ws.Range("D18").Value = 1.35
$ws.Range("E18").Value = 1.35
$ws.Range("D19").Value = 0.04
$ws.Range("E19").Value = 4.25
$ws.Range("C21").ClearContents()
$ws.Range("D22").Value = 5.78
$ws.Range("E22").Value = 4.39
$ws.Range("D23").Value = 2.84
$ws.Range("E23").Value = 2.95
$ws.Range("D24").Value = 4.25
$ws.Range("E24").Value = 0.5600000000000001
$ws.Range("D25").Value = 1.72
$ws.Range("E25").Value = 3.19
$ws.Range("D26").Value = 3.12
$ws.Range("E26").Value = 3.19
$ws.Range("D27").Value = 1.91
$ws.Range("D28").Value = 10.73
$ws.Range("E28").Value = 5.25
$ws.Range("D29").Value = 1.49
$ws.Range("E29").Value = 1.65
$ws.Range("D30").Value = 1.24
$ws.Range("E30").Value = 1.59
$ws.Range("C31").ClearContents()
$ws.Range("D32").Value = 2
$ws.Range("E32").Value = 0.75
$ws.Range("D33").Value = 2.12
$ws.Range("E33").Value = 0.34
$ws.Range("D34").Value = 6.11
$ws.Range("E34").Value = 1.85
$ws.Range("D35").Value = 4.42
$ws.Range("E35").Value = 1.49
$ws.Range("C36").ClearContents()
$ws.Range("D37").Value = 1.91
$ws.Range("E37").Value = 1.09
$ws.Range("D38").Value = 2.73
$ws.Range("E38").Value = 2.59
$ws.Range("C39").ClearContents()
$ws.Range("D40").Value = 1.99
$ws.Range("E40").Value = 1.99
$ws.Range("D41").Value = 3.72
$ws.Range("E41").Value = 3.15
$ws.Range("D42").Value = 0.33
$ws.Range("E42").Value = 7.3
$ws.Range("D43").Value = 5.39
$ws.Range("E43").Value = 14.95
$ws.Range("D44").Value = 26.45
$ws.Range("E44").Value = 20.9
$ws.Range("D45").Value = 2.96
$ws.Range("E45").Value = 1.32
$ws.Range("D46").Value = 23.92
$ws.Range("E46").Value = 1.55
$ws.Range("C47").ClearContents()
$ws.Range("D48").Value = 2.53
$ws.Range("E48").Value = 1.99
$ws.Range("D50").Value = 3.48
$ws.Range("E50").Value = 1.1
$ws.Range("C51").ClearContents()
$ws.Range("D52").Value = 11.08
$ws.Range("E52").Value = 4.5
$ws.Range("D53").Value = 7.14
$ws.Range("E53").Value = 1.75
$ws.Range("D54").Value = 23.79
$ws.Range("E54").Value = 5.99
$ws.Range("D55").Value = 14.96
$ws.Range("E55").Value = 5.99
$ws.Range("D56").Value = 12.5
$ws.Range("E56").Value = 3.9
$ws.Range("D57").Value = 0.67
$ws.Range("E57").Value = 1.38
$ws.Range("D58").Value = 4.43
$ws.Range("E58").Value = 3.45
$ws.Range("D59").Value = 1.54
$ws.Range("E59").Value = 1.43
$ws.Range("D60").Value = 2.98
$ws.Range("D61").Value = 10.95
$ws.Range("E61").Value = 2.29
$ws.Range("D62").Value = 3.98
$ws.Range("E62").Value = 2.95
$ws.Range("D63").Value = 5.7
$ws.Range("E63").Value = 3.99
$ws.Range("D64").Value = 4.42
$ws.Range("E64").Value = 5.95
$ws.Range("D65").Value = 14.56
$ws.Range("E65").Value = 13.99
$ws.Range("D66").Value = 14.56
$ws.Range("E66").Value = 10.1
$ws.Range("D67").Value = 4.54
$ws.Range("E67").Value = 6.45
$ws.Range("D68").Value = 1.05
$ws.Range("E68").Value = 1.24
$ws.Range("D69").Value = 1.05
$ws.Range("E69").Value = 1.24
$ws.Range("D70").Value = 1.24
$ws.Range("D71").Value = 1.15
$ws.Range("E71").Value = 1.19
$ws.Range("D72").Value = 3.36
$ws.Range("E72").Value = 3.46
$ws.Range("E74").Value = 1.35
$ws.Range("D75").Value = 3.07
$ws.Range("E75").Value = 2.75
$ws.Range("C76").ClearContents()
$ws.Range("D77").Value = 5.32
$ws.Range("E77").Value = 9.9
$ws.Range("D78").Value = 14.99
$ws.Range("E78").Value = 16.95
$ws.Range("D79").Value = 2.96
$ws.Range("E79").Value = 0.77
$ws.Range("D81").Value = 14.2
$ws.Range("E81").Value = 3.39
$ws.Range("D82").Value = 2.2
$ws.Range("E82").Value = 2.75
$ws.Range("D83").Value = 1.82
$ws.Range("E83").Value = 1.89
$ws.Range("D84").Value = 5.23
$ws.Range("E84").Value = 1.99
$ws.Range("D85").Value = 12.5
$ws.Range("E85").Value = 3.49
$ws.Range("C86").ClearContents()
$ws.Range("D88").Value = 1.82
$ws.Range("D89").Value = 10.24
$ws.Range("D90").Value = 11.75
$ws.Range("D91").Value = 7.2
$ws.Range("D93").Value = 2.9
$ws.Range("D94").Value = 20.91
$ws.Range("D95").Value = 16.33
$ws.Range("D96").Value = 0.12
$ws.Range("D97").Value = 0.14
$ws.Range("D98").Value = 0.5600000000000001
$ws.Range("D99").Value = 0.37
$ws.Range("D100").Value = 4.38
$ws.Range("D101").Value = 1.92
$ws.Range("D102").Value = 18.17
$ws.Range("D103").Value = 14.95
$ws.Range("D104").Value = 14.54
$ws.Range("D105").Value = 4.18
$ws.Range("D107").Value = 4.96
$ws.Range("D108").Value = 3.62
$ws.Range("D109").Value = 1.49
$ws.Range("D110").Value = 15.64
$ws.Range("D111").Value = 2.88
$ws.Range("D112").Value = 7.78
$ws.Range("D113").Value = 3.25
$ws.Range("D114").Value = 12.34
$ws.Range("D115").Value = 4.9
$ws.Range("C116").ClearContents()
$ws.Range("D117").Value = 13.03
$ws.Range("C118").ClearContents()
$ws.Range("D120").Value = 32.25
$ws.Range("C121").ClearContents()
$ws.Range("D122").Value = 8.539999999999999
$ws.Range("D123").Value = 13.83
$ws.Range("D124").Value = 0.06
$ws.Range("D125").Value = 2.96
$ws.Range("D127").Value = 2.88
$ws.Range("D128").Value = 4.86
$ws.Range("D129").Value = 3.94
$ws.Range("D130").Value = 6.9
$ws.Range("D131").Value = 7.96
$ws.Range("D132").Value = 6.09
$ws.Range("D133").Value = 1.32
$ws.Range("C134").ClearContents()
$ws.Range("D135").Value = 0.8
